$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Suite A (row 2) is no longer running
$ws.Range("C2").Value = "N"

# Suite E (row 6) and Suite F (row 7) are now running
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"

# Update the active selection to C7
$ws.Activate()
$ws.Range("C7").Select()
